$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 20.41270024389404
$ws.Range("C2").Value = 17.26868271090596
$ws.Range("D2").Value = 6.953944436805049
$ws.Range("E2").Value = 12.71111770342551
$ws.Range("F2").Value = 43.92638751237442
$ws.Range("I2").Value = 27.50587781547371
$ws.Range("J2").Value = 10.2843771402909
$ws.Range("N2").Value = 18.53984072357976
$ws.Range("B3").Value = 19.90059342218827
$ws.Range("C3").Value = 16.7600630216544
$ws.Range("D3").Value = 6.947847286123034
$ws.Range("E3").Value = 12.68434590529646
$ws.Range("F3").Value = 43.72777172053077
$ws.Range("I3").Value = 27.52309948137511
$ws.Range("J3").Value = 10.29380399114021
$ws.Range("N3").Value = 18.61783209767465
$ws.Range("B4").Value = 19.5854085791531
$ws.Range("C4").Value = 16.44545893365884
$ws.Range("D4").Value = 6.944803949421012
$ws.Range("E4").Value = 12.67057955769826
$ws.Range("F4").Value = 43.61938210404844
$ws.Range("I4").Value = 27.54139620639259
$ws.Range("J4").Value = 10.30152207838112
$ws.Range("N4").Value = 18.66772759342776
$ws.Range("B5").Value = 19.45698718252594
$ws.Range("C5").Value = 16.31689344150409
$ws.Range("D5").Value = 6.943741062016802
$ws.Range("E5").Value = 12.66564429437004
$ws.Range("F5").Value = 43.5786444692613
$ws.Range("I5").Value = 27.55078441393842
$ws.Range("J5").Value = 10.30515175345555
$ws.Range("N5").Value = 18.68856731357365
$ws.Range("B6").Value = 19.43567025350214
$ws.Range("C6").Value = 16.29553001385213
$ws.Range("D6").Value = 6.943575312388314
$ws.Range("E6").Value = 12.66486561524064
$ws.Range("F6").Value = 43.57208783411578
$ws.Range("I6").Value = 27.55245970965683
$ws.Range("J6").Value = 10.30578369402324
$ws.Range("N6").Value = 18.69205840309172
$ws.Range("B7").Value = 19.5836762872485
$ws.Range("C7").Value = 16.44372622136265
$ws.Range("D7").Value = 6.944788895539326
$ws.Range("E7").Value = 12.67051026412824
$ws.Range("F7").Value = 43.61881878140985
$ws.Range("I7").Value = 27.54151500977557
$ws.Range("J7").Value = 10.30156906894245
$ws.Range("N7").Value = 18.66800659039924
$ws.Range("B8").Value = 20.23640336403751
$ws.Range("C8").Value = 17.09391766427216
$ws.Range("D8").Value = 6.951697513356692
$ws.Range("E8").Value = 12.70133415287317
$ws.Range("F8").Value = 43.85510493491194
$ws.Range("I8").Value = 27.5102074714732
$ws.Range("J8").Value = 10.28722650177269
$ws.Range("N8").Value = 18.56631641237116
$ws.Range("B9").Value = 21.50176199423399
$ws.Range("C9").Value = 18.34163096372713
$ws.Range("D9").Value = 6.970748427323626
$ws.Range("E9").Value = 12.78282992754774
$ws.Range("F9").Value = 44.42484994014148
$ws.Range("I9").Value = 27.51047921114985
$ws.Range("J9").Value = 10.27444637018375
$ws.Range("N9").Value = 18.38275217405645
$ws.Range("B10").Value = 22.41155779875957
$ws.Range("C10").Value = 19.23062506491067
$ws.Range("D10").Value = 6.988025867660996
$ws.Range("E10").Value = 12.85530515753786
$ws.Range("F10").Value = 44.90632413087247
$ws.Range("I10").Value = 27.54873848264532
$ws.Range("J10").Value = 10.27444983223716
$ws.Range("N10").Value = 18.25742626349269
$ws.Range("B11").Value = 22.81904576916436
$ws.Range("C11").Value = 19.62699384937143
$ws.Range("D11").Value = 6.996581307703768
$ws.Range("E11").Value = 12.8909508854902
$ws.Range("F11").Value = 45.13848956576585
$ws.Range("I11").Value = 27.5744833857989
$ws.Range("J11").Value = 10.27649630325423
$ws.Range("N11").Value = 18.20245691061833
$ws.Range("B12").Value = 22.97227563811027
$ws.Range("C12").Value = 19.77578244515111
$ws.Range("D12").Value = 6.999919545511657
$ws.Range("E12").Value = 12.90482775805113
$ws.Range("F12").Value = 45.22824295537355
$ws.Range("I12").Value = 27.58543571784471
$ws.Range("J12").Value = 10.27756543703257
$ws.Range("N12").Value = 18.18193310800984
$ws.Range("B13").Value = 22.93932514299504
$ws.Range("C13").Value = 19.74379857610452
$ws.Range("D13").Value = 6.999196242228936
$ws.Range("E13").Value = 12.90182239162294
$ws.Range("F13").Value = 45.2088321133887
$ws.Range("I13").Value = 27.58302337128492
$ws.Range("J13").Value = 10.27732209588487
$ws.Range("N13").Value = 18.18634032150371
$ws.Range("B14").Value = 22.83167453372331
$ws.Range("C14").Value = 19.63926180056928
$ws.Range("D14").Value = 6.996853982151582
$ws.Range("E14").Value = 12.89208499392509
$ws.Range("F14").Value = 45.14583712083395
$ws.Range("I14").Value = 27.57536030696749
$ws.Range("J14").Value = 10.27657836557023
$ws.Range("N14").Value = 18.20076256490916
$ws.Range("B15").Value = 22.76559066457061
$ws.Range("C15").Value = 19.57505546142167
$ws.Range("D15").Value = 6.995432057065281
$ws.Range("E15").Value = 12.88616966825974
$ws.Range("F15").Value = 45.10748853077099
$ws.Range("I15").Value = 27.57082325818254
$ws.Range("J15").Value = 10.27616112121844
$ws.Range("N15").Value = 18.2096345676011
$ws.Range("B16").Value = 22.38478517165765
$ws.Range("C16").Value = 19.20454636521849
$ws.Range("D16").Value = 6.987480617656646
$ws.Range("E16").Value = 12.85302895870346
$ws.Range("F16").Value = 44.89141162515136
$ws.Range("I16").Value = 27.54722419169415
$ws.Range("J16").Value = 10.27435725342387
$ws.Range("N16").Value = 18.2610595816477
$ws.Range("B17").Value = 22.14942132504154
$ws.Range("C17").Value = 18.97507943351232
$ws.Range("D17").Value = 6.982779781204957
$ws.Range("E17").Value = 12.83337941316699
$ws.Range("F17").Value = 44.76218521014342
$ws.Range("I17").Value = 27.53488635312811
$ws.Range("J17").Value = 10.27377452511737
$ws.Range("N17").Value = 18.29312893217548
$ws.Range("B18").Value = 22.01345421336604
$ws.Range("C18").Value = 18.84234756145434
$ws.Range("D18").Value = 6.980141579060841
$ws.Range("E18").Value = 12.82232983657549
$ws.Range("F18").Value = 44.6890978333286
$ws.Range("I18").Value = 27.52857485462465
$ws.Range("J18").Value = 10.27363183290492
$ws.Range("N18").Value = 18.31176668390255
$ws.Range("B19").Value = 21.96732145644065
$ws.Range("C19").Value = 18.7972831269825
$ws.Range("D19").Value = 6.979259642413088
$ws.Range("E19").Value = 12.8186321581852
$ws.Range("F19").Value = 44.66456623745184
$ws.Range("I19").Value = 27.52657255601656
$ws.Range("J19").Value = 10.27361657158674
$ws.Range("N19").Value = 18.31811018895236
$ws.Range("B20").Value = 22.17453867651997
$ws.Range("C20").Value = 18.99958519880234
$ws.Range("D20").Value = 6.983273415766031
$ws.Range("E20").Value = 12.83544507004994
$ws.Range("F20").Value = 44.77581358758961
$ws.Range("I20").Value = 27.53611847040655
$ws.Range("J20").Value = 10.27381663542526
$ws.Range("N20").Value = 18.2896952044285
$ws.Range("B21").Value = 22.86332456185385
$ws.Range("C21").Value = 19.67000341083681
$ws.Range("D21").Value = 6.997539300936265
$ws.Range("E21").Value = 12.89493488074412
$ws.Range("F21").Value = 45.16429084743046
$ws.Range("I21").Value = 27.57757845263369
$ws.Range("J21").Value = 10.27678883301293
$ws.Range("N21").Value = 18.19651849422862
$ws.Range("B22").Value = 23.30714356487752
$ws.Range("C22").Value = 20.10047833547345
$ws.Range("D22").Value = 7.007436115206398
$ws.Range("E22").Value = 12.93601847040721
$ws.Range("F22").Value = 45.42886583607443
$ws.Range("I22").Value = 27.61168959050525
$ws.Range("J22").Value = 10.28044598763718
$ws.Range("N22").Value = 18.13732273180378
$ws.Range("B23").Value = 23.07089865601904
$ws.Range("C23").Value = 19.8714750108807
$ws.Range("D23").Value = 7.002102085687782
$ws.Range("E23").Value = 12.91389195877138
$ws.Range("F23").Value = 45.28669800056782
$ws.Range("I23").Value = 27.59284102768879
$ws.Range("J23").Value = 10.27833720346498
$ws.Range("N23").Value = 18.16876159255172
$ws.Range("B24").Value = 22.16318514459022
$ws.Range("C24").Value = 18.98850864824965
$ws.Range("D24").Value = 6.98305004295919
$ws.Range("E24").Value = 12.83451041610509
$ws.Range("F24").Value = 44.76964843591477
$ws.Range("I24").Value = 27.53555899571091
$ws.Range("J24").Value = 10.27379699831368
$ws.Range("N24").Value = 18.29124696761578
$ws.Range("B25").Value = 21.16218809327519
$ws.Range("C25").Value = 18.00821276725977
$ws.Range("D25").Value = 6.965012919905493
$ws.Range("E25").Value = 12.75855131019724
$ws.Range("F25").Value = 44.25951293928646
$ws.Range("I25").Value = 27.50375114788031
$ws.Range("J25").Value = 10.27625639133159
$ws.Range("N25").Value = 18.43072715036169
